## Trade #19 closed at 2026-02-16 22:53:57 - base_strategy UP +0.000%
## Append a new trade row (row 20) to both the "All Trades" and
## "base_strategy" worksheets, mirroring the existing row layout
## (Trade #, Date, Time, Strategy, Side, Entry Price, Exit Price, Status,
## P&L %, P&L $, Capital After, Entry Slippage, Exit Slippage, Confidence,
## Entry Reason, Exit Reason, Duration).

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 20

    # A: Trade #
    $ws.Cells.Item($row, 1).Value = 19

    # B: Date - force literal text (not an auto-converted date serial) by
    # writing it quote-prefixed, then resetting to the default "Normal"
    # style so no extra style index/quote-prefix flag sticks around.
    $ws.Cells.Item($row, 2).Value = "'2026-02-16"
    $ws.Cells.Item($row, 2).Style = "Normal"

    # C: Time (kept as plain text, matches source formatting)
    $ws.Cells.Item($row, 3).Value = "22:53:57"

    # D: Strategy
    $ws.Cells.Item($row, 4).Value = "base_strategy"

    # E: Side
    $ws.Cells.Item($row, 5).Value = "UP"

    # F: Entry Price
    $ws.Cells.Item($row, 6).Value = 49.999998

    # G: Exit Price (empty string, trade still open)
    $ws.Cells.Item($row, 7).Value = "'"
    $ws.Cells.Item($row, 7).Style = "Normal"

    # H: Status
    $ws.Cells.Item($row, 8).Value = "OPEN"

    # I: P&L %
    $ws.Cells.Item($row, 9).Value = 0

    # J: P&L $
    $ws.Cells.Item($row, 10).Value = 0

    # K: Capital After
    $ws.Cells.Item($row, 11).Value = 100

    # L: Entry Slippage (bps)
    $ws.Cells.Item($row, 12).Value = 0

    # M: Exit Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0

    # N: Confidence
    $ws.Cells.Item($row, 14).Value = 0.6

    # O: Entry Reason
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    # P: Exit Reason (empty string, trade still open)
    $ws.Cells.Item($row, 16).Value = "'"
    $ws.Cells.Item($row, 16).Style = "Normal"

    # Q: Duration (min)
    $ws.Cells.Item($row, 17).Value = 0
}
